$d = $word.ActiveDocument

$replacements = @(
    @{old="97×96=9312"; new="59×36=2124"},
    @{old="62×17=1054"; new="60×31=1860"},
    @{old="98×55=5390"; new="46×63=2898"},
    @{old="59×93=5487"; new="16×13=208"},
    @{old="17×81=1377"; new="56×62=3472"},
    @{old="19×36=684";  new="48×85=4080"},
    @{old="86×27=2322"; new="31×30=930"},
    @{old="17×69=1173"; new="84×76=6384"},
    @{old="60×39=2340"; new="16×38=608"},
    @{old="67×29=1943"; new="47×67=3149"},
    @{old="28×49=1372"; new="27×67=1809"},
    @{old="77×81=6237"; new="75×31=2325"},
    @{old="93×15=1395"; new="47×56=2632"},
    @{old="57×48=2736"; new="26×89=2314"},
    @{old="80×86=6880"; new="87×77=6699"},
    @{old="77×61=4697"; new="21×25=525"},
    @{old="22×19=418";  new="63×29=1827"},
    @{old="37×70=2590"; new="48×60=2880"},
    @{old="43×88=3784"; new="97×11=1067"},
    @{old="45×63=2835"; new="92×17=1564"},
    @{old="48×40=1920"; new="38×62=2356"},
    @{old="51×62=3162"; new="62×80=4960"},
    @{old="84×88=7392"; new="55×39=2145"},
    @{old="65×85=5525"; new="38×31=1178"},
    @{old="89×51=4539"; new="36×58=2088"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
